# Regenerate save_data to use K (strikeouts) instead of Strike# (total pitches
# thrown that were strikes), recalculated std/mean, and write the new s_vals
# into column G (header "K") for each game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column G ("K")
$newValues = @{
    2  = 11
    3  = 5
    4  = 2
    5  = 1
    6  = 8
    7  = 3
    8  = 10
    9  = 3
    10 = 4
    11 = 6
    12 = 2
    13 = 5
    14 = 0
    15 = 0
    16 = 1
    17 = 0
    18 = 2
    19 = 6
    20 = 1
    21 = 1
    22 = 8
    23 = 6
    24 = 2
    25 = 6
    26 = 12
    27 = 3
    28 = 3
    29 = 5
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}

$wb.Save()
